$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8983613252639771
$ws.Range("B1").Value = 1.571712970733643
$ws.Range("C1").Value = 6.446219921112061
$ws.Range("D1").Value = 2.970701932907104
$ws.Range("E1").Value = 1.548128128051758
